$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row at position 4, pushing existing rows 4:21 down to 5:22
$ws.Rows.Item(4).Insert()

# Populate the newly inserted row 4 with the new record
$ws.Cells.Item(4, 1).Value = 3
$ws.Cells.Item(4, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(4, 3).Value = "Coquimbo"
$ws.Cells.Item(4, 4).Value = 44592
$ws.Cells.Item(4, 5).Value = 5
$ws.Cells.Item(4, 6).Value = "Fruta"
$ws.Cells.Item(4, 7).Value = 100104
$ws.Cells.Item(4, 8).Value = "Frutos de pepita"
$ws.Cells.Item(4, 9).Value = 100104001
$ws.Cells.Item(4, 10).Value = "Granada"
$ws.Cells.Item(4, 11).Value = "Wonderfull"
$ws.Cells.Item(4, 12).Value = "Primera"
$ws.Cells.Item(4, 13).Value = 54
$ws.Cells.Item(4, 14).Value = 20000
$ws.Cells.Item(4, 15).Value = 20000
$ws.Cells.Item(4, 16).Value = 20000
$ws.Cells.Item(4, 17).Value = "$/caja 15 kilos empedrada"
$ws.Cells.Item(4, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(4, 19).Value = 1333
$ws.Cells.Item(4, 20).Value = 15
